$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 365.45834
$ws.Range("I92").Value = 392.72223
$ws.Range("J92").Value = 283.66666
$ws.Range("K92").Value = 392.72223
$ws.Range("L92").Value = 283.66666
$ws.Range("M92").Value = 855.2777699999999
$ws.Range("N92").Value = -2779.66666
$ws.Range("H117").Value = 99893
$ws.Range("J117").Value = 99893
$ws.Range("L117").Value = 99893
$ws.Range("N117").Value = -109071
$ws.Range("H123").Value = 70933.125
$ws.Range("J123").Value = 70933.125
$ws.Range("L123").Value = 70933.125
$ws.Range("N123").Value = -80733.125
$ws.Range("H135").Value = 1523.5454
$ws.Range("I135").Value = 973.2778
$ws.Range("K135").Value = 8759.5002
$ws.Range("M135").Value = -6224.5002
$ws.Range("H138").Value = 2111.5103
$ws.Range("I138").Value = 1863.027
$ws.Range("J138").Value = 2877.6667
$ws.Range("K138").Value = 5589.081
$ws.Range("L138").Value = 8633.000100000001
$ws.Range("M138").Value = -449.0810000000001
$ws.Range("N138").Value = -18913.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 611.0833
$ws.Range("I2").Value = 555.8889
$ws.Range("K2").Value = 555.8889
$ws.Range("M2").Value = -442.8889
$ws.Range("H32").Value = 7591.5
$ws.Range("I32").Value = 4839.878
$ws.Range("J32").Value = 14227.765
$ws.Range("K32").Value = 4839.878
$ws.Range("L32").Value = 14227.765
$ws.Range("M32").Value = -4552.878
$ws.Range("N32").Value = -14801.765
$ws.Range("H34").Value = 500000
$ws.Range("J34").Value = 500000
$ws.Range("L34").Value = 500000
$ws.Range("N34").Value = -500542
$ws.Range("H45").Value = 10427295
$ws.Range("I45").Value = 5363.4287
$ws.Range("K45").Value = 5363.4287
$ws.Range("M45").Value = -4986.4287
$ws.Range("H61").Value = 27345.275
$ws.Range("I61").Value = 2114.8438
$ws.Range("K61").Value = 2114.8438
$ws.Range("M61").Value = -1902.8438
$ws.Range("H63").Value = 2196
$ws.Range("I63").Value = 2196
$ws.Range("K63").Value = 2196
$ws.Range("M63").Value = -1510
$ws.Range("H66").Value = 2196
$ws.Range("I66").Value = 2196
$ws.Range("K66").Value = 10980
$ws.Range("M66").Value = -7548
$ws.Range("H116").Value = 611.0833
$ws.Range("I116").Value = 555.8889
$ws.Range("K116").Value = 555.8889
$ws.Range("M116").Value = 1738.1111
$ws.Range("H118").Value = 43569.145
$ws.Range("J118").Value = 43569.145
$ws.Range("L118").Value = 43569.145
$ws.Range("N118").Value = -46883.145
$ws.Range("H132").Value = 2152.4285
$ws.Range("I132").Value = 1649
$ws.Range("K132").Value = 4947
$ws.Range("M132").Value = -2417
$ws.Range("H136").Value = 27345.275
$ws.Range("I136").Value = 2114.8438
$ws.Range("K136").Value = 6344.5314
$ws.Range("M136").Value = -3794.5314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 611.0833
$ws.Range("I3").Value = 555.8889
$ws.Range("K3").Value = 555.8889
$ws.Range("M3").Value = -441.8889
$ws.Range("H52").Value = 99985
$ws.Range("J52").Value = 99985
$ws.Range("L52").Value = 99985
$ws.Range("N52").Value = -100511
$ws.Range("H55").Value = 36994.6
$ws.Range("J55").Value = 36994.6
$ws.Range("L55").Value = 36994.6
$ws.Range("N55").Value = -37540.6
$ws.Range("H114").Value = 89989
$ws.Range("J114").Value = 89989
$ws.Range("L114").Value = 89989
$ws.Range("N114").Value = -98667
$ws.Range("H115").Value = 82999.25
$ws.Range("J115").Value = 89998.664
$ws.Range("L115").Value = 89998.664
$ws.Range("N115").Value = -93132.664
$ws.Range("H121").Value = 99985
$ws.Range("J121").Value = 99985
$ws.Range("L121").Value = 99985
$ws.Range("N121").Value = -103479
$ws.Range("H127").Value = 47468.5
$ws.Range("J127").Value = 47468.5
$ws.Range("L127").Value = 47468.5
$ws.Range("N127").Value = -57388.5
$ws.Range("H134").Value = 2061.4722
$ws.Range("I134").Value = 1263.5217
$ws.Range("K134").Value = 3790.5651
$ws.Range("M134").Value = -1255.5651
$ws.Range("H138").Value = 89996
$ws.Range("J138").Value = 89996
$ws.Range("L138").Value = 89996
$ws.Range("N138").Value = -100276
$ws.Range("H140").Value = 43506.605
$ws.Range("J140").Value = 43557.5
$ws.Range("L140").Value = 43557.5
$ws.Range("N140").Value = -53917.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 32249.285
$ws.Range("J9").Value = 32249.285
$ws.Range("L9").Value = 32249.285
$ws.Range("N9").Value = -32585.285
$ws.Range("H16").Value = 1437.0834
$ws.Range("I16").Value = 1295.091
$ws.Range("K16").Value = 1295.091
$ws.Range("M16").Value = -1008.091
$ws.Range("H31").Value = 1871.9524
$ws.Range("I31").Value = 1555.6595
$ws.Range("K31").Value = 1555.6595
$ws.Range("M31").Value = -1260.6595
$ws.Range("H34").Value = 1871.9524
$ws.Range("I34").Value = 1555.6595
$ws.Range("K34").Value = 1555.6595
$ws.Range("M34").Value = -1353.6595
$ws.Range("H113").Value = 1437.0834
$ws.Range("I113").Value = 1295.091
$ws.Range("K113").Value = 1295.091
$ws.Range("M113").Value = 874.9090000000001
$ws.Range("H114").Value = 48700
$ws.Range("J114").Value = 48700
$ws.Range("L114").Value = 48700
$ws.Range("N114").Value = -57378
$ws.Range("H117").Value = 42641.145
$ws.Range("J117").Value = 42641.145
$ws.Range("L117").Value = 42641.145
$ws.Range("N117").Value = -51819.145
$ws.Range("H118").Value = 57997.715
$ws.Range("J118").Value = 57997.715
$ws.Range("L118").Value = 57997.715
$ws.Range("N118").Value = -61311.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 12902.8
$ws.Range("J132").Value = 12902.8
$ws.Range("L132").Value = 116125.2
$ws.Range("N132").Value = -121185.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 55040.09
$ws.Range("J108").Value = 55040.09
$ws.Range("L108").Value = 55040.09
$ws.Range("N108").Value = -62720.09
$ws.Range("H109").Value = 78408.625
$ws.Range("J109").Value = 78408.625
$ws.Range("L109").Value = 78408.625
$ws.Range("N109").Value = -80488.625
$ws.Range("H110").Value = 74800.164
$ws.Range("J110").Value = 74800.164
$ws.Range("L110").Value = 74800.164
$ws.Range("N110").Value = -82980.164
$ws.Range("H132").Value = 8636.182000000001
$ws.Range("I132").Value = 11874.5
$ws.Range("J132").Value = 7916.5557
$ws.Range("K132").Value = 35623.5
$ws.Range("L132").Value = 23749.6671
$ws.Range("M132").Value = -33093.5
$ws.Range("N132").Value = -28809.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 692.5625
$ws.Range("I16").Value = 692.5625
$ws.Range("K16").Value = 692.5625
$ws.Range("M16").Value = -522.5625
$ws.Range("H100").Value = 10460.137
$ws.Range("I100").Value = 11417.2
$ws.Range("J100").Value = 8409.286
$ws.Range("K100").Value = 11417.2
$ws.Range("L100").Value = 8409.286
$ws.Range("M100").Value = -10876.2
$ws.Range("N100").Value = -9491.286
$ws.Range("H121").Value = 42856
$ws.Range("J121").Value = 42856
$ws.Range("L121").Value = 42856
$ws.Range("N121").Value = -46350
$ws.Range("H133").Value = 86116
$ws.Range("J133").Value = 86116
$ws.Range("L133").Value = 86116
$ws.Range("N133").Value = -91176
$ws.Range("H136").Value = 2173.3157
$ws.Range("I136").Value = 1992.7142
$ws.Range("J136").Value = 2679
$ws.Range("K136").Value = 5978.142599999999
$ws.Range("L136").Value = 8037
$ws.Range("M136").Value = -3428.142599999999
$ws.Range("N136").Value = -13137

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 40848.4
$ws.Range("J121").Value = 40848.4
$ws.Range("L121").Value = 40848.4
$ws.Range("N121").Value = -44342.4
$ws.Range("H132").Value = 1012684.3
$ws.Range("I132").Value = 1236.4062
$ws.Range("J132").Value = 3955078.2
$ws.Range("K132").Value = 3709.2186
$ws.Range("L132").Value = 11865234.6
$ws.Range("M132").Value = -1179.2186
$ws.Range("N132").Value = -11870294.6
